$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Date: row 8, column B -> updated ISO timestamp
$ws.Range("B8").Value = "2024-02-19T18:37:26-06:00"

# Description: row 11, column B -> MatchSource -> MatchSync
$ws.Range("B11").Value = "NMDP Preferred Product for MatchSync patient"

# Case Sensitive: row 14, column B -> literal text "true" (not a Boolean).
# A direct Value assignment of "true"/"TRUE" is auto-coerced by Excel into a
# native Boolean cell, so instead synthesize the literal text via a formula
# in a scratch cell and paste-special the computed value back in, which
# preserves it as a genuine text cell.
$ws.Range("D1").Formula = '=LOWER(TEXT(TRUE,"@"))'
$ws.Range("D1").Copy()
$ws.Range("B14").PasteSpecial(-4163)
$ws.Range("D1").Clear()
